$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-46 (File_Name, Cell, Num_Inclusions)
$data = @(
    ,@("1K_PFF_ADAMTS19_01", 0, 30)
    ,@("1K_PFF_ADAMTS19_02", 0, 9)
    ,@("1K_PFF_ADAMTS19_02", 1, 9)
    ,@("1K_PFF_ADAMTS19_03", 0, 0)
    ,@("1K_PFF_ADAMTS19_03", 1, 17)
    ,@("1K_PFF_ADAMTS19_04", 0, 0)
    ,@("1K_PFF_ADAMTS19_04", 1, 2)
    ,@("1K_PFF_ADAMTS19_05", 0, 59)
    ,@("1K_PFF_ADAMTS19_05", 1, 11)
    ,@("1K_PFF_ADAMTS19_06", 0, 0)
    ,@("1K_PFF_ADAMTS19_06", 1, 0)
    ,@("1K_PFF_ADAMTS19_06", 2, 0)
    ,@("1K_PFF_ADAMTS19_06", 3, 12)
    ,@("1K_PFF_ADAMTS19_06", 4, 0)
    ,@("1K_PFF_ADAMTS19_06", 5, 0)
    ,@("1K_PFF_ADAMTS19_07", 0, 20)
    ,@("1K_PFF_ADAMTS19_08", 0, 0)
    ,@("1K_PFF_ADAMTS19_09", 0, 0)
    ,@("1K_PFF_ADAMTS19_10", 0, 0)
    ,@("1K_PFF_ADAMTS19_10", 1, 0)
    ,@("1K_PFF_ADAMTS19_10", 2, 0)
    ,@("1K_PFF_SCR_01", 0, 118)
    ,@("1K_PFF_SCR_01", 1, 0)
    ,@("1K_PFF_SCR_02", 0, 26)
    ,@("1K_PFF_SCR_03", 0, 0)
    ,@("1K_PFF_SCR_04", 0, 0)
    ,@("1K_PFF_SCR_05", 0, 0)
    ,@("1K_PFF_SCR_06", 0, 0)
    ,@("1K_PFF_SCR_07", 0, 0)
    ,@("1K_PFF_SCR_08", 0, 1)
    ,@("1K_PFF_SCR_09", 0, 0)
    ,@("1K_PFF_SCR_10", 0, 0)
    ,@("1K_PFF_TAX1BP1_01", 0, 0)
    ,@("1K_PFF_TAX1BP1_01", 1, 3)
    ,@("1K_PFF_TAX1BP1_02", 0, 18)
    ,@("1K_PFF_TAX1BP1_03", 0, 0)
    ,@("1K_PFF_TAX1BP1_04", 0, 0)
    ,@("1K_PFF_TAX1BP1_05", 0, 0)
    ,@("1K_PFF_TAX1BP1_06", 0, 0)
    ,@("1K_PFF_TAX1BP1_06", 1, 0)
    ,@("1K_PFF_TAX1BP1_07", 0, 0)
    ,@("1K_PFF_TAX1BP1_07", 1, 0)
    ,@("1K_PFF_TAX1BP1_08", 0, 0)
    ,@("1K_PFF_TAX1BP1_09", 0, 0)
    ,@("1K_PFF_TAX1BP1_10", 0, 0)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Remove now-unused rows 47-65 (sheet shrank from 65 to 46 rows)
$ws.Range("A47:C65").ClearContents()

